$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 187.4
$ws.Range("I15").Value = 187.4
$ws.Range("K15").Value = 562.2
$ws.Range("M15").Value = -393.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 263.42426
$ws.Range("I33").Value = 246.14815
$ws.Range("K33").Value = 246.14815
$ws.Range("M33").Value = -17.14814999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 260.2
$ws.Range("I5").Value = 266.33334
$ws.Range("J5").Value = 251
$ws.Range("K5").Value = 266.33334
$ws.Range("L5").Value = 251
$ws.Range("M5").Value = -154.33334
$ws.Range("N5").Value = -475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5007495
$ws.Range("I74").Value = 6431259
$ws.Range("J74").Value = 71780
$ws.Range("K74").Value = 6431259
$ws.Range("L74").Value = 71780
$ws.Range("M74").Value = -6430385
$ws.Range("N74").Value = -73528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5007495
$ws.Range("I77").Value = 6431259
$ws.Range("J77").Value = 71780
$ws.Range("K77").Value = 32156295
$ws.Range("L77").Value = 358900
$ws.Range("M77").Value = -32151927
$ws.Range("N77").Value = -367636

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 50429
$ws.Range("J128").Value = 50429
$ws.Range("L128").Value = 50429
$ws.Range("N128").Value = -60389

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 260.2
$ws.Range("I4").Value = 266.33334
$ws.Range("J4").Value = 251
$ws.Range("K4").Value = 266.33334
$ws.Range("L4").Value = 251
$ws.Range("M4").Value = -151.33334
$ws.Range("N4").Value = -481

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1001.2857
$ws.Range("I94").Value = 802.25
$ws.Range("J94").Value = 1266.6666
$ws.Range("K94").Value = 802.25
$ws.Range("L94").Value = 1266.6666
$ws.Range("M94").Value = -351.25
$ws.Range("N94").Value = -2168.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11114179
$ws.Range("I134").Value = 2970.9412
$ws.Range("K134").Value = 8912.8236
$ws.Range("M134").Value = -6377.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 82293.2
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 82293.2
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 82293.2
$ws.Range("N31").Value = -82883.2
$ws.Range("M31").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 82293.2
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 82293.2
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 82293.2
$ws.Range("N34").Value = -82697.2
$ws.Range("M34").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12822053
$ws.Range("I58").Value = 22728600
$ws.Range("J58").Value = 1814.8529
$ws.Range("K58").Value = 22728600
$ws.Range("L58").Value = 1814.8529
$ws.Range("M58").Value = -22728397
$ws.Range("N58").Value = -2220.8529

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 29761.805
$ws.Range("I132").Value = 1697.6923
$ws.Range("J132").Value = 102728.5
$ws.Range("K132").Value = 5093.0769
$ws.Range("L132").Value = 308185.5
$ws.Range("M132").Value = -2563.0769
$ws.Range("N132").Value = -313245.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 12822053
$ws.Range("I136").Value = 22728600
$ws.Range("J136").Value = 1814.8529
$ws.Range("K136").Value = 68185800
$ws.Range("L136").Value = 5444.5587
$ws.Range("M136").Value = -68183250
$ws.Range("N136").Value = -10544.5587

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 498.5
$ws.Range("I108").Value = 498.5
$ws.Range("K108").Value = 1495.5
$ws.Range("M108").Value = 1384.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 963.5
$ws.Range("I110").Value = 963.5
$ws.Range("K110").Value = 2890.5
$ws.Range("M110").Value = 1199.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 879
$ws.Range("I113").Value = 546.6667
$ws.Range("J113").Value = 989.7778
$ws.Range("K113").Value = 1640.0001
$ws.Range("L113").Value = 2969.3334
$ws.Range("M113").Value = 529.9999
$ws.Range("N113").Value = -7309.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 3116.923
$ws.Range("I120").Value = 2376.6667
$ws.Range("K120").Value = 7130.000100000001
$ws.Range("M120").Value = -2292.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 937.2432
$ws.Range("I122").Value = 247.4
$ws.Range("J122").Value = 1192.7407
$ws.Range("K122").Value = 2226.6
$ws.Range("L122").Value = 10734.6663
$ws.Range("M122").Value = 223.4000000000001
$ws.Range("N122").Value = -15634.6663

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 857.1385
$ws.Range("I131").Value = 464.45456
$ws.Range("J131").Value = 937.12964
$ws.Range("K131").Value = 1393.36368
$ws.Range("L131").Value = 2811.38892
$ws.Range("M131").Value = 3646.63632
$ws.Range("N131").Value = -12891.38892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3919
$ws.Range("I141").Value = 3919
$ws.Range("K141").Value = 11757
$ws.Range("M141").Value = -6577

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1834.8334
$ws.Range("I122").Value = 1506.0294
$ws.Range("K122").Value = 4518.0882
$ws.Range("M122").Value = -2068.0882

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 74302.71000000001
$ws.Range("I132").Value = 45628.87
$ws.Range("J132").Value = 206202.4
$ws.Range("K132").Value = 136886.61
$ws.Range("L132").Value = 618607.2
$ws.Range("M132").Value = -134356.61
$ws.Range("N132").Value = -623667.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 24928.572
$ws.Range("J135").Value = 24928.572
$ws.Range("L135").Value = 24928.572
$ws.Range("N135").Value = -35068.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1444356.5
$ws.Range("I46").Value = 4329468.5
$ws.Range("J46").Value = 1800.4286
$ws.Range("K46").Value = 4329468.5
$ws.Range("L46").Value = 1800.4286
$ws.Range("M46").Value = -4329280.5
$ws.Range("N46").Value = -2176.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1946.4736
$ws.Range("I61").Value = 1836.7693
$ws.Range("K61").Value = 1836.7693
$ws.Range("M61").Value = -1634.7693

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1946.4736
$ws.Range("I113").Value = 1836.7693
$ws.Range("K113").Value = 1836.7693
$ws.Range("M113").Value = 333.2307000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 42177.926
$ws.Range("J134").Value = 42177.926
$ws.Range("L134").Value = 42177.926
$ws.Range("N134").Value = -52317.926

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 35577.21
$ws.Range("J139").Value = 35577.21
$ws.Range("L139").Value = 35577.21
$ws.Range("N139").Value = -45857.21

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 34271.71
$ws.Range("I136").Value = 21900.75
$ws.Range("J136").Value = 76686.42999999999
$ws.Range("K136").Value = 65702.25
$ws.Range("L136").Value = 230059.29
$ws.Range("M136").Value = -63152.25
$ws.Range("N136").Value = -235159.29
